# "Opponent photo. Add friends dialog button." —
# appends five new backlog items to the "Бэклог задач" sheet (sheet2):
#   B58/C58 .. B61/C61 : task text + timestamp (same date style as the
#                        existing rows above, i.e. row 57's C-column format)
#   B62                : a final task text with no timestamp yet
# and leaves the selection on the next free row (A62), ready for the next
# entry to be typed in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 58
$ws.Range("B58").Value = "вынести ID_* профайлера в отедельный файл"
$ws.Range("C58").Value = 42040.138888888891
$ws.Range("C57").Copy() | Out-Null
$ws.Range("C58").PasteSpecial(-4122) | Out-Null

# Row 59
$ws.Range("B59").Value = "GUIDom not selectable by default"
$ws.Range("C59").Value = 42040.140277777777
$ws.Range("C57").Copy() | Out-Null
$ws.Range("C59").PasteSpecial(-4122) | Out-Null

# Row 60
$ws.Range("B60").Value = "Рейтинг не верно отображается(сортировка?!)"
$ws.Range("C60").Value = 42040.140277777777
$ws.Range("C57").Copy() | Out-Null
$ws.Range("C60").PasteSpecial(-4122) | Out-Null

# Row 61
$ws.Range("B61").Value = "Рейтинг перекрывает онлайн и очки"
$ws.Range("C61").Value = 42040.140277777777
$ws.Range("C57").Copy() | Out-Null
$ws.Range("C61").PasteSpecial(-4122) | Out-Null

# Row 62 - text only, no timestamp yet
$ws.Range("B62").Value = "ElementPhoto - effect old photo like-a"

# Leave the cursor on the next empty row, like the sheet was left mid-edit
$ws.Range("A62").Select() | Out-Null
